# literature-review.xlsx : proposal and start introduction
# Adds a new "large language model cybersecurity" keyword-search block
# (rows 70-80) to the "Paper selection" sheet, mirroring the layout of
# the existing keyword blocks above it (Keywords/Count/Titles/related
# flags), then repositions the sheet view to where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70: first row of the new keyword block (Keyword + Count + first title) ---
# Shared strings must be created in this exact order (C before A) so the
# new sharedStrings table matches the source ordering.
$ws.Range("C70").Value = "CySecBERT: A Domain-Adapted Language Model for the Cybersecurity Domain"
$ws.Range("A70").Value = "large language model cybersecurity"
$ws.Range("B70").Value = 119000
$ws.Range("D70").Value = $true

# --- Rows 71-80: remaining papers found for this keyword search ---
$ws.Range("C71").Value = "Securebert: A domain-specific language model for cybersecurity"
$ws.Range("D71").Value = $true

$ws.Range("C72").Value = "Cybert: Cybersecurity claim classification by fine-tuning the bert language model"
$ws.Range("D72").Value = $false

$ws.Range("C73").Value = "Cybert: Contextualized embeddings for the cybersecurity domain"
$ws.Range("D73").Value = $false

$ws.Range("C74").Value = "Generative AI and Large Language Modeling in Cybersecurity"
$ws.Range("D74").Value = $true

$ws.Range("C75").Value = "Assessing Cybersecurity Vulnerabilities in Code Large Language Models"
$ws.Range("D75").Value = $false

$ws.Range("C76").Value = "Large Language Models in Cybersecurity: State-of-the-Art"
$ws.Range("D76").Value = $true

$ws.Range("C77").Value = "CyberSecEval 2: A Wide-Ranging Cybersecurity Evaluation Suite for Large Language Models"
$ws.Range("D77").Value = $true

$ws.Range("C78").Value = "What are the latest cybersecurity trends? a case study grounded in language models"
$ws.Range("D78").Value = $true

$ws.Range("C79").Value = "Using large language models for cybersecurity capture-the-flag challenges and certification questions"
$ws.Range("D79").Value = $false

$ws.Range("C80").Value = "AI, ML, AND LARGE LANGUAGE MODELS IN CYBERSECURITY"
$ws.Range("D80").Value = $true
$ws.Range("E80").Value = $true

# --- Reposition the view / selection to where editing left off ---
$ws.Range("C74").Select() | Out-Null
